{"js": "// Replace the date line and each division-problem cell's text, in document\n// order, with the corresponding new value. All 26 original text values are\n// unique in this document, so an exact-match search + full-text replace is\n// safe and unambiguous.\nconst replacements = [\n  [\"2025-03-07 Friday\", \"2025-03-08 Saturday\"],\n  [\"50\u00f73=16, 2\", \"51\u00f75=10, 1\"],\n  [\"99\u00f73=33, 0\", \"23\u00f78=2, 7\"],\n  [\"89\u00f74=22, 1\", \"78\u00f79=8, 6\"],\n  [\"42\u00f72=21, 0\", \"82\u00f73=27, 1\"],\n  [\"65\u00f78=8, 1\", \"34\u00f76=5, 4\"],\n  [\"83\u00f73=27, 2\", \"87\u00f78=10, 7\"],\n  [\"54\u00f76=9, 0\", \"40\u00f74=10, 0\"],\n  [\"60\u00f75=12, 0\", \"71\u00f77=10, 1\"],\n  [\"97\u00f79=10, 7\", \"18\u00f78=2, 2\"],\n  [\"50\u00f72=25, 0\", \"77\u00f79=8, 5\"],\n  [\"81\u00f78=10, 1\", \"50\u00f72=25, 0\"],\n  [\"32\u00f72=16, 0\", \"92\u00f75=18, 2\"],\n  [\"61\u00f75=12, 1\", \"70\u00f75=14, 0\"],\n  [\"19\u00f78=2, 3\", \"22\u00f79=2, 4\"],\n  [\"58\u00f77=8, 2\", \"60\u00f79=6, 6\"],\n  [\"33\u00f74=8, 1\", \"36\u00f77=5, 1\"],\n  [\"19\u00f79=2, 1\", \"48\u00f74=12, 0\"],\n  [\"56\u00f76=9, 2\", \"66\u00f73=22, 0\"],\n  [\"46\u00f73=15, 1\", \"76\u00f79=8, 4\"],\n  [\"25\u00f76=4, 1\", \"16\u00f74=4, 0\"],\n  [\"67\u00f74=16, 3\", \"67\u00f76=11, 1\"],\n  [\"27\u00f72=13, 1\", \"55\u00f76=9, 1\"],\n  [\"46\u00f72=23, 0\", \"62\u00f78=7, 6\"],\n  [\"47\u00f75=9, 2\", \"69\u00f73=23, 0\"],\n  [\"11\u00f77=1, 4\", \"42\u00f76=7, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Exactly one occurrence is expected for every value in this document.\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each division-problem cell's text, in document\n# order, with the corresponding new value. All 26 original text values are\n# unique in this document, so an exact-match Find/Replace (MatchCase +\n# MatchWholeWord, replace first occurrence only) is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-03-07 Friday\", \"2025-03-08 Saturday\"),\n    @(\"50\u00f73=16, 2\", \"51\u00f75=10, 1\"),\n    @(\"99\u00f73=33, 0\", \"23\u00f78=2, 7\"),\n    @(\"89\u00f74=22, 1\", \"78\u00f79=8, 6\"),\n    @(\"42\u00f72=21, 0\", \"82\u00f73=27, 1\"),\n    @(\"65\u00f78=8, 1\", \"34\u00f76=5, 4\"),\n    @(\"83\u00f73=27, 2\", \"87\u00f78=10, 7\"),\n    @(\"54\u00f76=9, 0\", \"40\u00f74=10, 0\"),\n    @(\"60\u00f75=12, 0\", \"71\u00f77=10, 1\"),\n    @(\"97\u00f79=10, 7\", \"18\u00f78=2, 2\"),\n    @(\"50\u00f72=25, 0\", \"77\u00f79=8, 5\"),\n    @(\"81\u00f78=10, 1\", \"50\u00f72=25, 0\"),\n    @(\"32\u00f72=16, 0\", \"92\u00f75=18, 2\"),\n    @(\"61\u00f75=12, 1\", \"70\u00f75=14, 0\"),\n    @(\"19\u00f78=2, 3\", \"22\u00f79=2, 4\"),\n    @(\"58\u00f77=8, 2\", \"60\u00f79=6, 6\"),\n    @(\"33\u00f74=8, 1\", \"36\u00f77=5, 1\"),\n    @(\"19\u00f79=2, 1\", \"48\u00f74=12, 0\"),\n    @(\"56\u00f76=9, 2\", \"66\u00f73=22, 0\"),\n    @(\"46\u00f73=15, 1\", \"76\u00f79=8, 4\"),\n    @(\"25\u00f76=4, 1\", \"16\u00f74=4, 0\"),\n    @(\"67\u00f74=16, 3\", \"67\u00f76=11, 1\"),\n    @(\"27\u00f72=13, 1\", \"55\u00f76=9, 1\"),\n    @(\"46\u00f72=23, 0\", \"62\u00f78=7, 6\"),\n    @(\"47\u00f75=9, 2\", \"69\u00f73=23, 0\"),\n    @(\"11\u00f77=1, 4\", \"42\u00f76=7, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # Args: FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #       MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #       Format, ReplaceWith, Replace(2 = wdReplaceOne)\n    $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
